$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new items as boosts (item name in column C, item index in column D)
$ws.Range("C24").Value = "Cheese boost"
$ws.Range("D24").Value = 6

$ws.Range("C25").Value = "Bacon boost"
$ws.Range("D25").Value = 7

$ws.Range("C26").Value = "Garlic boost"
$ws.Range("D26").Value = 8

$ws.Range("C27").Value = "Basil boost"
$ws.Range("D27").Value = 9

# Update the view state to match (scrolled down, new selection)
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("D28").Select()
